$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Type-changing cells: copy number-format from a donor cell, set value, fix type ---
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C22").Formula = '="0"'
$ws.Range("C22").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("I26").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C26").Value = 1

$ws.Range("I26").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("F26").Value = 1

$ws.Range("C30").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("D30").Formula = '="0"'
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("E27").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("E30").Value = "***.*"

# --- Simple numeric value updates ---
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 53.333333333333
$ws.Range("I16").Value = 173
$ws.Range("J16").Value = 158
$ws.Range("K16").Value = 9.493670886075
$ws.Range("L16").Value = 42.97520661157
$ws.Range("M16").Value = -48.816568047337
$ws.Range("N16").Value = -84.995663486556
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = -23.076923076923
$ws.Range("I17").Value = 269
$ws.Range("J17").Value = 255
$ws.Range("K17").Value = 5.490196078431
$ws.Range("L17").Value = 36.548223350253
$ws.Range("M17").Value = 46.195652173913
$ws.Range("N17").Value = -49.436090225563
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 20
$ws.Range("G18").Value = 30
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 248
$ws.Range("J18").Value = 274
$ws.Range("K18").Value = -9.48905109489
$ws.Range("L18").Value = 24
$ws.Range("M18").Value = -40.811455847255
$ws.Range("N18").Value = -79.419087136929
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -7.692307692307
$ws.Range("F19").Value = 64
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = -3.030303030303
$ws.Range("I19").Value = 685
$ws.Range("J19").Value = 608
$ws.Range("K19").Value = 12.66447368421
$ws.Range("L19").Value = 50.219298245614
$ws.Range("M19").Value = 53.24384787472
$ws.Range("N19").Value = 40.368852459016
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 16.666666666666
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -21.052631578947
$ws.Range("I20").Value = 163
$ws.Range("J20").Value = 166
$ws.Range("K20").Value = -1.807228915662
$ws.Range("L20").Value = 19.85294117647
$ws.Range("M20").Value = 15.602836879432
$ws.Range("N20").Value = -80.337756332931
$ws.Range("C21").Value = 33
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -2.941176470588
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 156
$ws.Range("H21").Value = -8.974358974358
$ws.Range("I21").Value = 1554
$ws.Range("J21").Value = 1476
$ws.Range("K21").Value = 5.284552845528
$ws.Range("L21").Value = 37.400530503978
$ws.Range("M21").Value = 0.974658869395
$ws.Range("N21").Value = -63.529687866698
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 21
$ws.Range("K22").Value = -41.666666666666
$ws.Range("L22").Value = 5
$ws.Range("M22").Value = -47.5
$ws.Range("C23").Value = 6
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = 50
$ws.Range("F23").Value = 16
$ws.Range("G23").Value = 20
$ws.Range("H23").Value = -20
$ws.Range("I23").Value = 179
$ws.Range("J23").Value = 146
$ws.Range("K23").Value = 22.602739726027
$ws.Range("L23").Value = 20.134228187919
$ws.Range("M23").Value = 45.528455284552
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 81
$ws.Range("H24").Value = -33.333333333333
$ws.Range("I24").Value = 921
$ws.Range("J24").Value = 1068
$ws.Range("K24").Value = -13.76404494382
$ws.Range("L24").Value = 4.421768707482
$ws.Range("M24").Value = -16.877256317689
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 45
$ws.Range("G25").Value = 51
$ws.Range("H25").Value = -11.764705882352
$ws.Range("I25").Value = 442
$ws.Range("J25").Value = 451
$ws.Range("K25").Value = -1.995565410199
$ws.Range("L25").Value = 31.940298507462
$ws.Range("M25").Value = 0
$ws.Range("I26").Value = 20
$ws.Range("K26").Value = -4.761904761904
$ws.Range("L26").Value = -13.043478260869
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -42.857142857142
$ws.Range("I27").Value = 50
$ws.Range("K27").Value = 13.636363636363
$ws.Range("L27").Value = -20.63492063492
$ws.Range("N28").Value = -85.365853658536
$ws.Range("N29").Value = -88.157894736842
$ws.Range("F30").Value = 3
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 18
$ws.Range("K30").Value = -37.931034482758
$ws.Range("L30").Value = 50
